# Applies "msz - table part 1" changes to Tabelle1 (sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 updates -------------------------------------------------
# B8/C8 keep their existing text (no edit needed).
# H8 is removed entirely.
$ws.Range("H8").ClearContents()

# New "check" hint cells (written first so their shared-string entries
# are created before the "enter" cells below).
$ws.Range("E8").Value = "Product Page check for hint date with invalid format"
$ws.Range("G8").Value = "Product Page check for hint date with invalid value in past"

# D8 / F8 get new text content and keep the Text ("@") number format.
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "Product Page enter date with invalid format"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "Product Page enter date with invalid value in past"

# --- New row 9 -------------------------------------------------------
$ws.Range("A9").Value = "102_AutomobileInsurance_005_PriceOption_001_MandatoryFields"
$ws.Range("A10").Value = "102_AutomobileInsurance_005_PriceOption_002_FieldHintsAndErrors"

$ws.Range("B9").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageVehicleData"
$ws.Range("C9").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageInsurantData"
$ws.Range("D9").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageProductData"
$ws.Range("E9").Value = "Goto price option page"
$ws.Range("F9").Value = "Price option page check for open mandatory field"
$ws.Range("H9").Value = "Price option page check for filled mandatory field"
$ws.Range("G9").Value = "Select Ultimate"

# --- New row 10 -------------------------------------------------------
$ws.Range("B10").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageVehicleData"
$ws.Range("C10").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageInsurantData"
$ws.Range("D10").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageProductData"
$ws.Range("E10").Value = "Goto price option page"

# --- Cosmetic sheet/window adjustments -------------------------------
$ws.Columns("E:E").ColumnWidth = 43
$ws.Range("A11:XFD14").Select()
